# Apply "add minimal information gain to database - version 2"
$wb = $excel.ActiveWorkbook

# Work on the "Options" sheet (second sheet in the workbook)
$ws = $wb.Worksheets.Item("Options")

# Add new row: key/value pair, matching the existing Title/Number participants rows
$ws.Range("A3").Value = "Minimal information gain"
$ws.Range("B3").Value = 10

# Update selection on the Options sheet, as in the diff
$ws.Range("A4").Select()

# Make Options the active (selected) sheet/tab
$ws.Activate()
